# v 1.0.18 Resolucion incidencia
# Swap the two posting-line details between row 2 and row 3 (same customer /
# document, two line items: "WAGE ADV.GIVEN PERS." (debit) and
# "RECEIVABL.FROM PERS." (credit) were recorded against the wrong line -
# correct it by exchanging Item/Posting Key/CME/G-L Account/G-L Account
# Descr./Debit amount/Credit amount between the two rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: becomes the "RECEIVABL.FROM PERS." line ---
$ws.Range("Q2").Value = "1"
$ws.Range("R2").Value = "17"
$ws.Range("S2").Value = "1"
$ws.Range("T2").Value = "1350101001"
$ws.Range("U2").Value = "RECEIVABL.FROM PERS."
$ws.Range("V2").Value = 0
$ws.Range("X2").Value = 1625.25

# --- Row 3: becomes the "WAGE ADV.GIVEN PERS." line ---
$ws.Range("Q3").Value = "2"
$ws.Range("R3").Value = "09"
$ws.Range("S3").Value = "0"
$ws.Range("T3").Value = "1960101001"
$ws.Range("U3").Value = "WAGE ADV.GIVEN PERS."
$ws.Range("V3").Value = 1625.25
$ws.Range("X3").Value = 0
